# Updates the "Price" (D) and "Volume(1h)" (E) columns of the cryptos
# worksheet with refreshed values.
#
# Many of the "Price" strings look like plain numbers (e.g. "604.80",
# "5.04"), and Excel's automatic type detection would otherwise coerce
# them into numeric cells (losing the original text formatting / trailing
# zeros, e.g. "604.80" -> 604.79999999999995). To keep them as literal
# text - matching the source data, which stores every cell as a string -
# each value is written with a leading apostrophe (forces "text" entry)
# and the cell style is reset back to "Normal" right after so no stray
# quote-prefix/number-format style lingers on the cell.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''70.998.78'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '''  +4.87%  '
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = '''2.620.56'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '''  +5.52%  '
$ws.Range("E3").Style = "Normal"
$ws.Range("E4").Value = '''  +0.00%  '
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = '''604.80'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '''  +3.02%  '
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = '''181.03'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '''  +3.71%  '
$ws.Range("E6").Style = "Normal"
$ws.Range("E7").Value = '''  -0.05%  '
$ws.Range("E7").Style = "Normal"
$ws.Range("E8").Value = '''  +2.10%  '
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = '''2.619.24'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '''  +5.51%  '
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = '''0.164'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '''  +13.68%  '
$ws.Range("E10").Style = "Normal"
$ws.Range("E11").Value = '''  +0.39%  '
$ws.Range("E11").Style = "Normal"
$ws.Range("E12").Value = '''  +4.64%  '
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = '''5.04'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '''  +1.84%  '
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = '''3.069.92'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '''  +4.50%  '
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = '''26.71'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '''  +5.94%  '
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = '''0.0000182'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '''  +7.76%  '
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = '''71.029.85'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '''  +4.80%  '
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = '''2.610.93'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '''  +4.60%  '
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = '''381.87'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '''  +10.27%  '
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = '''7.93'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '''  +7.66%  '
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = '''11.46'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '''  +6.37%  '
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = '''4.17'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '''  +1.87%  '
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = '''71.96'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '''  +1.83%  '
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = '''4.44'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '''  +6.42%  '
$ws.Range("E24").Style = "Normal"
$ws.Range("E25").Value = '''  +0.05%  '
$ws.Range("E25").Style = "Normal"
$ws.Range("E26").Value = '''  +12.02%  '
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = '''9.68'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '''  +9.79%  '
$ws.Range("E27").Style = "Normal"
$ws.Range("D28").Value = '''2.750.88'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '''  +5.34%  '
$ws.Range("E28").Style = "Normal"
$ws.Range("D29").Value = '''0.998'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '''  -0.04%  '
$ws.Range("E29").Style = "Normal"
$ws.Range("D30").Value = '''0.0₃0948'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '''  +6.65%  '
$ws.Range("E30").Style = "Normal"
$ws.Range("D31").Value = '''526.61'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '''  +5.73%  '
$ws.Range("E31").Style = "Normal"
$ws.Range("D32").Value = '''8.03'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '''  +4.40%  '
$ws.Range("E32").Style = "Normal"
$ws.Range("E33").Value = '''  +7.15%  '
$ws.Range("E33").Style = "Normal"
$ws.Range("E34").Value = '''  +4.45%  '
$ws.Range("E34").Style = "Normal"
$ws.Range("E35").Value = '''  -0.01%  '
$ws.Range("E35").Style = "Normal"
$ws.Range("D36").Value = '''164.68'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '''  +0.07%  '
$ws.Range("E36").Style = "Normal"
$ws.Range("E37").Value = '''  +0.16%  '
$ws.Range("E37").Style = "Normal"
$ws.Range("D38").Value = '''1.92'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '''  +12.25%  '
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").Value = '''19.14'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '''  +4.95%  '
$ws.Range("E39").Style = "Normal"
$ws.Range("E40").Value = '''  +1.71%  '
$ws.Range("E40").Style = "Normal"
$ws.Range("E41").Value = '''  +6.28%  '
$ws.Range("E41").Style = "Normal"
$ws.Range("E42").Value = '''  +0.10%  '
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = '''5.04'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '''  +6.00%  '
$ws.Range("E43").Style = "Normal"
$ws.Range("E44").Value = '''  +9.06%  '
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = '''0.331'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '''  +2.46%  '
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = '''40.17'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '''  +3.93%  '
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = '''153.75'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '''  +3.92%  '
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = '''3.65'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '''  +4.06%  '
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = '''0.0₆0270'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '''  +7.23%  '
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = '''0.533'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '''  +4.67%  '
$ws.Range("E50").Style = "Normal"
$ws.Range("E51").Value = '''  +7.32%  '
$ws.Range("E51").Style = "Normal"
